$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells remain text, matching the original inline-string
# representation, so plain-looking decimal values are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.946.03'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '2.305.20'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '305.59'
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("D6").Value = '97.49'
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("D7").Value = '0.511'
$ws.Range("E7").Value = '  -1.34%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.504'
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("D10").Value = '35.72'
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("D11").Value = '0.0796'
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '18.31'
$ws.Range("E12").Value = '  +1.90%  '
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("E14").Value = '  -1.46%  '
$ws.Range("D15").Value = '2.661.45'
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").Value = '2.299.03'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '0.784'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").Value = '42.871.91'
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("D19").Value = '13.14'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").Value = '0.0₃0905'
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").Value = '6.06'
$ws.Range("D22").Value = '67.64'
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("D23").Value = '236.94'
$ws.Range("E23").Value = '  -0.52%  '
$ws.Range("D24").Value = '2.17'
$ws.Range("E24").Value = '  -1.81%  '
$ws.Range("E25").Value = '  +2.56%  '
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("D28").Value = '25.52'
$ws.Range("E28").Value = '  +1.35%  '
$ws.Range("D29").Value = '167.45'
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("D30").Value = '2.07'
$ws.Range("E30").Value = '  +1.50%  '
$ws.Range("D31").Value = '9.10'
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").Value = '33.09'
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("D34").Value = '4.83'
$ws.Range("E34").Value = '  +1.92%  '
$ws.Range("E35").Value = '  -2.25%  '
$ws.Range("D36").Value = '17.45'
$ws.Range("E36").Value = '  -5.11%  '
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").Value = '1.76'
$ws.Range("E40").Value = '  -1.64%  '
$ws.Range("E41").Value = '  -0.96%  '
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("D43").Value = '2.020.71'
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("D44").Value = '0.0281'
$ws.Range("E44").Value = '  -2.35%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '18.11'
$ws.Range("E45").Value = '  +3.85%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '10.03'
$ws.Range("E46").Value = '  -2.10%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '2.12'
$ws.Range("E47").Value = '  -3.08%  '
$ws.Range("E48").Value = '  -1.76%  '
$ws.Range("E49").Value = '  +6.06%  '
$ws.Range("D50").Value = '54.06'
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").Value = '2.527.97'
$ws.Range("E51").Value = '  -0.34%  '
